$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 28: update title and link
$ws.Range("D28").Value = "강화학습 논문 정리 15편 : Accelerating Interactive Human-like Manipulation Learning with GPU-based Simulation and High-quality Demonstrations (IEEE-RAS 2022)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/231"

# Row 32: update title and link
$ws.Range("D32").Value = "[디자인 패턴] Singleton 싱글턴 패턴"
$ws.Range("E32").Value = "https://dodonam.tistory.com/467"

# Row 51: update title and link
$ws.Range("D51").Value = "[nextjs] tailwindcss 사용할 때 globals.css에 있어야 하는 코드"
$ws.Range("E51").Value = "https://bskyvision.com/entry/nextjs-tailwindcss-%EC%82%AC%EC%9A%A9%ED%95%A0-%EB%95%8C-globalscss%EC%97%90-%EC%9E%88%EC%96%B4%EC%95%BC-%ED%95%98%EB%8A%94-%EC%BD%94%EB%93%9C"
